$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("D7").Value = 31.80303075250784
$ws.Range("F7").Value = $false

# Row 8
$ws.Range("D8").Value = 37.14197833583991

# Row 12
$ws.Range("D12").Value = 14.87776628947827

# Row 18
$ws.Range("F18").Value = $false

# Row 23
$ws.Range("D23").Value = 8.009948565874655

# Row 29
$ws.Range("E29").Value = 20
$ws.Range("F29").Value = $false

# Row 34
$ws.Range("D34").Value = 11.10533207739218
